$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6976
$ws1.Range("F4").Value = 114
$ws1.Range("F5").Value = 163
$ws1.Range("F8").Value = 594

# Sheet "全部类型" (All Types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6976
$ws4.Range("F5").Value = 114
$ws4.Range("F6").Value = 163
$ws4.Range("F10").Value = 594
